# Auto-generated edit script applying the cryptos.xlsx diff (updated crypto prices/volumes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.NumberFormat = "General"
}

Set-TextValue $ws.Range("D2") "61.691.39"
Set-TextValue $ws.Range("E2") "  +2.54%  "
Set-TextValue $ws.Range("D3") "3.392.94"
Set-TextValue $ws.Range("E3") "  +1.61%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D5") "575.78"
Set-TextValue $ws.Range("E5") "  +1.79%  "
Set-TextValue $ws.Range("D6") "136.14"
Set-TextValue $ws.Range("E6") "  +3.54%  "
Set-TextValue $ws.Range("E7") "  -0.03%  "
Set-TextValue $ws.Range("B8") "LidoStakedEther"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextValue $ws.Range("D8") "3.390.97"
Set-TextValue $ws.Range("E8") "  +1.56%  "
Set-TextValue $ws.Range("B9") "XRP"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D9") "0.475"
Set-TextValue $ws.Range("E9") "  +0.26%  "
Set-TextValue $ws.Range("B10") "Toncoin"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D10") "7.45"
Set-TextValue $ws.Range("E10") "  -0.14%  "
Set-TextValue $ws.Range("B11") "Dogecoin"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D11") "0.126"
Set-TextValue $ws.Range("E11") "  +6.07%  "
Set-TextValue $ws.Range("B12") "Cardano"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D12") "0.391"
Set-TextValue $ws.Range("E12") "  +3.83%  "
Set-TextValue $ws.Range("B13") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D13") "3.974.34"
Set-TextValue $ws.Range("E13") "  +1.66%  "
Set-TextValue $ws.Range("B14") "TRON"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D14") "0.122"
Set-TextValue $ws.Range("E14") "  +2.47%  "
Set-TextValue $ws.Range("B15") "ShibaInu"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D15") "0.0000176"
Set-TextValue $ws.Range("E15") "  +3.59%  "
Set-TextValue $ws.Range("B16") "WrappedEther"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D16") "3.400.17"
Set-TextValue $ws.Range("E16") "  +1.88%  "
Set-TextValue $ws.Range("B17") "Avalanche"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D17") "25.33"
Set-TextValue $ws.Range("E17") "  +2.03%  "
Set-TextValue $ws.Range("B18") "WrappedBTC"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D18") "61.772.74"
Set-TextValue $ws.Range("E18") "  +2.50%  "
Set-TextValue $ws.Range("B19") "Chainlink"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D19") "14.19"
Set-TextValue $ws.Range("E19") "  +4.70%  "
Set-TextValue $ws.Range("B20") "Polkadot"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D20") "5.84"
Set-TextValue $ws.Range("E20") "  +2.36%  "
Set-TextValue $ws.Range("B21") "Uniswap"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D21") "9.44"
Set-TextValue $ws.Range("E21") "  +2.06%  "
Set-TextValue $ws.Range("B22") "BitcoinCash"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D22") "386.75"
Set-TextValue $ws.Range("E22") "  +8.84%  "
Set-TextValue $ws.Range("B23") "Polygon"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D23") "0.567"
Set-TextValue $ws.Range("E23") "  +1.39%  "
Set-TextValue $ws.Range("B24") "WrappedeETH"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D24") "3.539.75"
Set-TextValue $ws.Range("E24") "  +2.04%  "
Set-TextValue $ws.Range("B25") "Dai"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D25") "1.00"
Set-TextValue $ws.Range("E25") "  -0.03%  "
Set-TextValue $ws.Range("B26") "PEPE"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D26") "0.0000127"
Set-TextValue $ws.Range("E26") "  +13.89%  "
Set-TextValue $ws.Range("B27") "Litecoin"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D27") "71.19"
Set-TextValue $ws.Range("E27") "  +2.57%  "
Set-TextValue $ws.Range("B28") "RenderToken"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D28") "7.64"
Set-TextValue $ws.Range("E28") "  +2.26%  "
Set-TextValue $ws.Range("B29") "Fetch.AI"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D29") "1.58"
Set-TextValue $ws.Range("E29") "  -3.11%  "
Set-TextValue $ws.Range("B30") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  +0.17%  "
Set-TextValue $ws.Range("B31") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D31") "8.26"
Set-TextValue $ws.Range("E31") "  +3.57%  "
Set-TextValue $ws.Range("B32") "Kaspa"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D32") "0.161"
Set-TextValue $ws.Range("E32") "  +4.47%  "
Set-TextValue $ws.Range("B33") "PancakeSwap"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D33") "2.17"
Set-TextValue $ws.Range("E33") "  +1.94%  "
Set-TextValue $ws.Range("B34") "USDe"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D34") "1.00"
Set-TextValue $ws.Range("E34") "  +0.04%  "
Set-TextValue $ws.Range("B35") "RenzoRestakedETH"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Range("D35") "3.424.30"
Set-TextValue $ws.Range("E35") "  +1.62%  "
Set-TextValue $ws.Range("B36") "EthereumClassic"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D36") "23.48"
Set-TextValue $ws.Range("E36") "  +2.31%  "
Set-TextValue $ws.Range("B37") "NEARProtocol"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D37") "5.40"
Set-TextValue $ws.Range("E37") "  +0.75%  "
Set-TextValue $ws.Range("B38") "Aptos"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D38") "6.95"
Set-TextValue $ws.Range("E38") "  +0.87%  "
Set-TextValue $ws.Range("B39") "ImmutableX"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D39") "1.55"
Set-TextValue $ws.Range("E39") "  +2.33%  "
Set-TextValue $ws.Range("B40") "Monero"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D40") "163.31"
Set-TextValue $ws.Range("E40") "  +2.74%  "
Set-TextValue $ws.Range("B41") "Hedera"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D41") "0.0785"
Set-TextValue $ws.Range("E41") "  +1.58%  "
Set-TextValue $ws.Range("B42") "Stacks"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "1.78"
Set-TextValue $ws.Range("E42") "  +12.55%  "
Set-TextValue $ws.Range("B43") "Mantle"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.783"
Set-TextValue $ws.Range("E43") "  +4.43%  "
Set-TextValue $ws.Range("B44") "FirstDigitalUSD"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D44") "1.00"
Set-TextValue $ws.Range("E44") "  +0.06%  "
Set-TextValue $ws.Range("B45") "ONDO"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws.Range("D45") "1.23"
Set-TextValue $ws.Range("E45") "  +3.30%  "
Set-TextValue $ws.Range("B46") "Filecoin"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D46") "4.44"
Set-TextValue $ws.Range("E46") "  +1.11%  "
Set-TextValue $ws.Range("B47") "OKB"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D47") "41.69"
Set-TextValue $ws.Range("E47") "  +1.93%  "
Set-TextValue $ws.Range("B48") "EnergySwap"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D48") "24.60"
Set-TextValue $ws.Range("E48") "  +4.37%  "
Set-TextValue $ws.Range("B49") "Cosmos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D49") "6.94"
Set-TextValue $ws.Range("E49") "  +1.92%  "
Set-TextValue $ws.Range("B50") "InjectiveProtocol"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D50") "23.25"
Set-TextValue $ws.Range("E50") "  +3.43%  "
Set-TextValue $ws.Range("B51") "Maker"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D51") "2.357.91"
Set-TextValue $ws.Range("E51") "  +8.08%  "
